# Append rows 43-55 to the worksheet, replicating the source export.
# Every field is written as literal text, matching the original rows 2-42
# (all stored as text, even the numeric- and boolean-looking values).
# A leading "'" forces Excel to store the literal text rather than
# auto-converting "true"/"false" to booleans, or long digit strings to
# floating point numbers (which would lose precision).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43
$ws.Cells.Item(43, 1).Value = "'624042747312230"
$ws.Cells.Item(43, 2).Value = "'true"
$ws.Cells.Item(43, 3).Value = "'0"
$ws.Cells.Item(43, 4).Value = "'0"
$ws.Cells.Item(43, 5).Value = "'0"
$ws.Cells.Item(43, 6).Value = "'237663744490"
$ws.Cells.Item(43, 7).Value = "'true"
$ws.Cells.Item(43, 8).Value = "'3587660809968478"
$ws.Cells.Item(43, 9).Value = "'None"
$ws.Cells.Item(43, 10).Value = "'true"
$ws.Cells.Item(43, 11).Value = "'237660002052"
$ws.Cells.Item(43, 12).Value = "'10.124.148.4"
$ws.Cells.Item(43, 13).Value = "'None"
$ws.Cells.Item(43, 14).Value = "'None"
$ws.Cells.Item(43, 15).Value = "'None"
$ws.Cells.Item(43, 16).Value = "'None"
$ws.Cells.Item(43, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(43, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(43, 19).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(43, 20).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(43, 21).Value = "'DOMS02:KNOWN SUBSCRIBER;result:ok;"

# Row 44
$ws.Cells.Item(44, 1).Value = "'624042747312230"
$ws.Cells.Item(44, 2).Value = "'true"
$ws.Cells.Item(44, 3).Value = "'0"
$ws.Cells.Item(44, 4).Value = "'0"
$ws.Cells.Item(44, 5).Value = "'0"
$ws.Cells.Item(44, 6).Value = "'237663744490"
$ws.Cells.Item(44, 7).Value = "'true"
$ws.Cells.Item(44, 8).Value = "'3587660809968478"
$ws.Cells.Item(44, 9).Value = "'None"
$ws.Cells.Item(44, 10).Value = "'true"
$ws.Cells.Item(44, 11).Value = "'237660002052"
$ws.Cells.Item(44, 12).Value = "'10.124.148.4"
$ws.Cells.Item(44, 13).Value = "'None"
$ws.Cells.Item(44, 14).Value = "'None"
$ws.Cells.Item(44, 15).Value = "'None"
$ws.Cells.Item(44, 16).Value = "'None"
$ws.Cells.Item(44, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(44, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(44, 19).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(44, 20).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(44, 21).Value = "'DOMS02:KNOWN SUBSCRIBER;result:ok;"

# Row 45
$ws.Cells.Item(45, 1).Value = "'624042747827156"
$ws.Cells.Item(45, 2).Value = "'true"
$ws.Cells.Item(45, 3).Value = "'0"
$ws.Cells.Item(45, 4).Value = "'2"
$ws.Cells.Item(45, 5).Value = "'0"
$ws.Cells.Item(45, 6).Value = "'237669595858"
$ws.Cells.Item(45, 7).Value = "'true"
$ws.Cells.Item(45, 8).Value = "'8630780379935655"
$ws.Cells.Item(45, 9).Value = "'None"
$ws.Cells.Item(45, 10).Value = "'true"
$ws.Cells.Item(45, 11).Value = "'237660002051"
$ws.Cells.Item(45, 12).Value = "'10.124.140.1"
$ws.Cells.Item(45, 13).Value = "'None"
$ws.Cells.Item(45, 14).Value = "'None"
$ws.Cells.Item(45, 15).Value = "'None"
$ws.Cells.Item(45, 16).Value = "'None"
$ws.Cells.Item(45, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(45, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(45, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(45, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(45, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;odbic:Barring ic solved;"

# Row 46
$ws.Cells.Item(46, 1).Value = "'624042747827156"
$ws.Cells.Item(46, 2).Value = "'true"
$ws.Cells.Item(46, 3).Value = "'0"
$ws.Cells.Item(46, 4).Value = "'2"
$ws.Cells.Item(46, 5).Value = "'0"
$ws.Cells.Item(46, 6).Value = "'237669595858"
$ws.Cells.Item(46, 7).Value = "'true"
$ws.Cells.Item(46, 8).Value = "'8630780379935655"
$ws.Cells.Item(46, 9).Value = "'None"
$ws.Cells.Item(46, 10).Value = "'true"
$ws.Cells.Item(46, 11).Value = "'237660002051"
$ws.Cells.Item(46, 12).Value = "'10.124.140.1"
$ws.Cells.Item(46, 13).Value = "'None"
$ws.Cells.Item(46, 14).Value = "'None"
$ws.Cells.Item(46, 15).Value = "'None"
$ws.Cells.Item(46, 16).Value = "'None"
$ws.Cells.Item(46, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(46, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(46, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(46, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(46, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;odbic:Barring ic solved;"

# Row 47
$ws.Cells.Item(47, 1).Value = "'624042747827156"
$ws.Cells.Item(47, 2).Value = "'true"
$ws.Cells.Item(47, 3).Value = "'0"
$ws.Cells.Item(47, 4).Value = "'0"
$ws.Cells.Item(47, 5).Value = "'0"
$ws.Cells.Item(47, 6).Value = "'237669595858"
$ws.Cells.Item(47, 7).Value = "'true"
$ws.Cells.Item(47, 8).Value = "'8630780379935655"
$ws.Cells.Item(47, 9).Value = "'None"
$ws.Cells.Item(47, 10).Value = "'true"
$ws.Cells.Item(47, 11).Value = "'237660002051"
$ws.Cells.Item(47, 12).Value = "'10.124.140.1"
$ws.Cells.Item(47, 13).Value = "'None"
$ws.Cells.Item(47, 14).Value = "'None"
$ws.Cells.Item(47, 15).Value = "'None"
$ws.Cells.Item(47, 16).Value = "'None"
$ws.Cells.Item(47, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(47, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(47, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(47, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(47, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;result:ok;"

# Row 48
$ws.Cells.Item(48, 1).Value = "'624042747827156"
$ws.Cells.Item(48, 2).Value = "'true"
$ws.Cells.Item(48, 3).Value = "'0"
$ws.Cells.Item(48, 4).Value = "'0"
$ws.Cells.Item(48, 5).Value = "'0"
$ws.Cells.Item(48, 6).Value = "'237669595858"
$ws.Cells.Item(48, 7).Value = "'true"
$ws.Cells.Item(48, 8).Value = "'8630780379935655"
$ws.Cells.Item(48, 9).Value = "'None"
$ws.Cells.Item(48, 10).Value = "'true"
$ws.Cells.Item(48, 11).Value = "'237660002051"
$ws.Cells.Item(48, 12).Value = "'10.124.140.1"
$ws.Cells.Item(48, 13).Value = "'None"
$ws.Cells.Item(48, 14).Value = "'None"
$ws.Cells.Item(48, 15).Value = "'None"
$ws.Cells.Item(48, 16).Value = "'None"
$ws.Cells.Item(48, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(48, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(48, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(48, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(48, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;result:ok;"

# Row 49
$ws.Cells.Item(49, 1).Value = "'624042747827156"
$ws.Cells.Item(49, 2).Value = "'true"
$ws.Cells.Item(49, 3).Value = "'0"
$ws.Cells.Item(49, 4).Value = "'2"
$ws.Cells.Item(49, 5).Value = "'0"
$ws.Cells.Item(49, 6).Value = "'237669595858"
$ws.Cells.Item(49, 7).Value = "'true"
$ws.Cells.Item(49, 8).Value = "'8630780379935655"
$ws.Cells.Item(49, 9).Value = "'None"
$ws.Cells.Item(49, 10).Value = "'true"
$ws.Cells.Item(49, 11).Value = "'237660002051"
$ws.Cells.Item(49, 12).Value = "'10.124.140.1"
$ws.Cells.Item(49, 13).Value = "'None"
$ws.Cells.Item(49, 14).Value = "'None"
$ws.Cells.Item(49, 15).Value = "'None"
$ws.Cells.Item(49, 16).Value = "'None"
$ws.Cells.Item(49, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(49, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(49, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(49, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(49, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;odbic:Barring ic solved;"

# Row 50
$ws.Cells.Item(50, 1).Value = "'624042747827156"
$ws.Cells.Item(50, 2).Value = "'true"
$ws.Cells.Item(50, 3).Value = "'1"
$ws.Cells.Item(50, 4).Value = "'2"
$ws.Cells.Item(50, 5).Value = "'0"
$ws.Cells.Item(50, 6).Value = "'237669595858"
$ws.Cells.Item(50, 7).Value = "'true"
$ws.Cells.Item(50, 8).Value = "'8630780379935655"
$ws.Cells.Item(50, 9).Value = "'None"
$ws.Cells.Item(50, 10).Value = "'true"
$ws.Cells.Item(50, 11).Value = "'237660002051"
$ws.Cells.Item(50, 12).Value = "'10.124.140.1"
$ws.Cells.Item(50, 13).Value = "'None"
$ws.Cells.Item(50, 14).Value = "'None"
$ws.Cells.Item(50, 15).Value = "'None"
$ws.Cells.Item(50, 16).Value = "'None"
$ws.Cells.Item(50, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(50, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(50, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(50, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(50, 21).Value = "'odboc:Barring oc solved;DOMS01:KNOWN SUBSCRIBER;odbic:Barring ic solved;"

# Row 51
$ws.Cells.Item(51, 1).Value = "'624042747827156"
$ws.Cells.Item(51, 2).Value = "'true"
$ws.Cells.Item(51, 3).Value = "'0"
$ws.Cells.Item(51, 4).Value = "'0"
$ws.Cells.Item(51, 5).Value = "'0"
$ws.Cells.Item(51, 6).Value = "'237669595858"
$ws.Cells.Item(51, 7).Value = "'true"
$ws.Cells.Item(51, 8).Value = "'8630780379935655"
$ws.Cells.Item(51, 9).Value = "'None"
$ws.Cells.Item(51, 10).Value = "'true"
$ws.Cells.Item(51, 11).Value = "'237660002051"
$ws.Cells.Item(51, 12).Value = "'10.124.140.1"
$ws.Cells.Item(51, 13).Value = "'None"
$ws.Cells.Item(51, 14).Value = "'None"
$ws.Cells.Item(51, 15).Value = "'None"
$ws.Cells.Item(51, 16).Value = "'None"
$ws.Cells.Item(51, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(51, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(51, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(51, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(51, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;result:ok;"

# Row 52
$ws.Cells.Item(52, 1).Value = "'624042747827156"
$ws.Cells.Item(52, 2).Value = "'true"
$ws.Cells.Item(52, 3).Value = "'1"
$ws.Cells.Item(52, 4).Value = "'2"
$ws.Cells.Item(52, 5).Value = "'0"
$ws.Cells.Item(52, 6).Value = "'237669595858"
$ws.Cells.Item(52, 7).Value = "'true"
$ws.Cells.Item(52, 8).Value = "'8630780379935655"
$ws.Cells.Item(52, 9).Value = "'None"
$ws.Cells.Item(52, 10).Value = "'true"
$ws.Cells.Item(52, 11).Value = "'237660002051"
$ws.Cells.Item(52, 12).Value = "'10.124.140.1"
$ws.Cells.Item(52, 13).Value = "'None"
$ws.Cells.Item(52, 14).Value = "'None"
$ws.Cells.Item(52, 15).Value = "'None"
$ws.Cells.Item(52, 16).Value = "'None"
$ws.Cells.Item(52, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(52, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(52, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(52, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(52, 21).Value = "'odboc:Barring oc solved;DOMS01:KNOWN SUBSCRIBER;odbic:Barring ic solved;"

# Row 53
$ws.Cells.Item(53, 1).Value = "'624042747827156"
$ws.Cells.Item(53, 2).Value = "'true"
$ws.Cells.Item(53, 3).Value = "'0"
$ws.Cells.Item(53, 4).Value = "'0"
$ws.Cells.Item(53, 5).Value = "'0"
$ws.Cells.Item(53, 6).Value = "'237669595858"
$ws.Cells.Item(53, 7).Value = "'true"
$ws.Cells.Item(53, 8).Value = "'8630780379935655"
$ws.Cells.Item(53, 9).Value = "'None"
$ws.Cells.Item(53, 10).Value = "'true"
$ws.Cells.Item(53, 11).Value = "'237660002051"
$ws.Cells.Item(53, 12).Value = "'10.124.140.1"
$ws.Cells.Item(53, 13).Value = "'None"
$ws.Cells.Item(53, 14).Value = "'None"
$ws.Cells.Item(53, 15).Value = "'None"
$ws.Cells.Item(53, 16).Value = "'None"
$ws.Cells.Item(53, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(53, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(53, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(53, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(53, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;result:ok;"

# Row 54
$ws.Cells.Item(54, 1).Value = "'624042747827156"
$ws.Cells.Item(54, 2).Value = "'true"
$ws.Cells.Item(54, 3).Value = "'0"
$ws.Cells.Item(54, 4).Value = "'0"
$ws.Cells.Item(54, 5).Value = "'0"
$ws.Cells.Item(54, 6).Value = "'237669595858"
$ws.Cells.Item(54, 7).Value = "'true"
$ws.Cells.Item(54, 8).Value = "'8630780379935655"
$ws.Cells.Item(54, 9).Value = "'None"
$ws.Cells.Item(54, 10).Value = "'true"
$ws.Cells.Item(54, 11).Value = "'237660002051"
$ws.Cells.Item(54, 12).Value = "'10.124.140.1"
$ws.Cells.Item(54, 13).Value = "'None"
$ws.Cells.Item(54, 14).Value = "'None"
$ws.Cells.Item(54, 15).Value = "'None"
$ws.Cells.Item(54, 16).Value = "'None"
$ws.Cells.Item(54, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(54, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(54, 19).Value = "'KNOWN SUBSCRIBER"
$ws.Cells.Item(54, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(54, 21).Value = "'DOMS01:KNOWN SUBSCRIBER;result:ok;"

# Row 55
$ws.Cells.Item(55, 1).Value = "'None"
$ws.Cells.Item(55, 2).Value = "'None"
$ws.Cells.Item(55, 3).Value = "'None"
$ws.Cells.Item(55, 4).Value = "'None"
$ws.Cells.Item(55, 5).Value = "'None"
$ws.Cells.Item(55, 6).Value = "'None"
$ws.Cells.Item(55, 7).Value = "'None"
$ws.Cells.Item(55, 8).Value = "'None"
$ws.Cells.Item(55, 9).Value = "'searchResponse:requestID=ee4874dc-cedf-4198-b96e-be2e9aef8cac, errorCode=32, errorMessage= error result (32); matchedDN = dc=MSISDN,DC=C-NTDB,entries:`nEND OF SEARCH ENTRIES.,"
$ws.Cells.Item(55, 10).Value = "'None"
$ws.Cells.Item(55, 11).Value = "'None"
$ws.Cells.Item(55, 12).Value = "'None"
$ws.Cells.Item(55, 13).Value = "'None"
$ws.Cells.Item(55, 14).Value = "'None"
$ws.Cells.Item(55, 15).Value = "'None"
$ws.Cells.Item(55, 16).Value = "'None"
$ws.Cells.Item(55, 17).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(55, 18).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(55, 19).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(55, 20).Value = "'UNKNOWN SUBSCRIBER"
$ws.Cells.Item(55, 21).Value = "'ldapResponse:Unknow Subscriber in HLR;"
